$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-13 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-14 Friday", 2) | Out-Null
$d.Content.Find.Execute("689×4=2756", $true, $false, $false, $false, $false, $true, 1, $false, "844×5=4220", 2) | Out-Null
$d.Content.Find.Execute("707×5=3535", $true, $false, $false, $false, $false, $true, 1, $false, "375×6=2250", 2) | Out-Null
$d.Content.Find.Execute("393×2=786", $true, $false, $false, $false, $false, $true, 1, $false, "234×8=1872", 2) | Out-Null
$d.Content.Find.Execute("502×5=2510", $true, $false, $false, $false, $false, $true, 1, $false, "876×6=5256", 2) | Out-Null
$d.Content.Find.Execute("403×4=1612", $true, $false, $false, $false, $false, $true, 1, $false, "984×2=1968", 2) | Out-Null
$d.Content.Find.Execute("639×6=3834", $true, $false, $false, $false, $false, $true, 1, $false, "957×5=4785", 2) | Out-Null
$d.Content.Find.Execute("680×2=1360", $true, $false, $false, $false, $false, $true, 1, $false, "114×2=228", 2) | Out-Null
$d.Content.Find.Execute("923×9=8307", $true, $false, $false, $false, $false, $true, 1, $false, "991×7=6937", 2) | Out-Null
$d.Content.Find.Execute("631×5=3155", $true, $false, $false, $false, $false, $true, 1, $false, "710×2=1420", 2) | Out-Null
$d.Content.Find.Execute("396×8=3168", $true, $false, $false, $false, $false, $true, 1, $false, "539×5=2695", 2) | Out-Null
$d.Content.Find.Execute("742×8=5936", $true, $false, $false, $false, $false, $true, 1, $false, "589×3=1767", 2) | Out-Null
$d.Content.Find.Execute("426×3=1278", $true, $false, $false, $false, $false, $true, 1, $false, "380×3=1140", 2) | Out-Null
$d.Content.Find.Execute("312×4=1248", $true, $false, $false, $false, $false, $true, 1, $false, "162×7=1134", 2) | Out-Null
$d.Content.Find.Execute("410×3=1230", $true, $false, $false, $false, $false, $true, 1, $false, "520×2=1040", 2) | Out-Null
$d.Content.Find.Execute("855×6=5130", $true, $false, $false, $false, $false, $true, 1, $false, "315×6=1890", 2) | Out-Null
$d.Content.Find.Execute("652×8=5216", $true, $false, $false, $false, $false, $true, 1, $false, "638×2=1276", 2) | Out-Null
$d.Content.Find.Execute("690×5=3450", $true, $false, $false, $false, $false, $true, 1, $false, "367×9=3303", 2) | Out-Null
$d.Content.Find.Execute("475×8=3800", $true, $false, $false, $false, $false, $true, 1, $false, "642×7=4494", 2) | Out-Null
$d.Content.Find.Execute("900×3=2700", $true, $false, $false, $false, $false, $true, 1, $false, "493×7=3451", 2) | Out-Null
$d.Content.Find.Execute("198×4=792", $true, $false, $false, $false, $false, $true, 1, $false, "490×9=4410", 2) | Out-Null
$d.Content.Find.Execute("267×3=801", $true, $false, $false, $false, $false, $true, 1, $false, "534×8=4272", 2) | Out-Null
$d.Content.Find.Execute("458×4=1832", $true, $false, $false, $false, $false, $true, 1, $false, "747×9=6723", 2) | Out-Null
$d.Content.Find.Execute("793×4=3172", $true, $false, $false, $false, $false, $true, 1, $false, "210×7=1470", 2) | Out-Null
$d.Content.Find.Execute("165×6=990", $true, $false, $false, $false, $false, $true, 1, $false, "499×9=4491", 2) | Out-Null
$d.Content.Find.Execute("444×2=888", $true, $false, $false, $false, $false, $true, 1, $false, "234×9=2106", 2) | Out-Null
